# Generate Report for Archive
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it appears
#    (Overview!E:F, zh-cn!C, de-de!C).
# 2) Narrow the corresponding "Status" columns to match the new, shorter text.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $found = $used.Find("Ready for handoff")
    if ($found -ne $null) {
        $firstAddress = $found.Address()
        while ($true) {
            $found.Value = "In Translation"
            $found = $used.FindNext($found)
            if ($found -eq $null -or $found.Address() -eq $firstAddress) {
                break
            }
        }
    }
}

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
